# Recompute Y_UTM (B), X_UTM (C), ZoneNumber (D) and ZoneLetter (E) for every
# data row, forcing UTM zone 48 / zone letter "T" instead of the previously
# auto-detected zone (1 / "W"). Easting/Northing are derived from the
# existing Latitude (I) / Longitude (J) columns using the standard UTM
# projection formulas (WGS84 ellipsoid).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wf = $excel.WorksheetFunction

$PI = 3.14159265358979323846

function Radians($deg) {
    return $deg * $PI / 180.0
}

# Forward UTM projection (lat/lon in degrees) forced onto zone 48, letter T.
function ToUtmZone48T($lat, $lon) {
    $K0 = 0.9996
    $E = 0.00669438
    $E_P2 = $E / (1 - $E)
    $R = 6378137.0

    $M1 = (1 - $E / 4 - 3 * $E * $E / 64 - 5 * $E * $E * $E / 256)
    $M2 = (3 * $E / 8 + 3 * $E * $E / 32 + 45 * $E * $E * $E / 1024)
    $M3 = (15 * $E * $E / 256 + 45 * $E * $E * $E / 1024)
    $M4 = (35 * $E * $E * $E / 3072)

    $lat_rad = Radians $lat
    $lat_sin = $wf.Sin($lat_rad)
    $lat_cos = $wf.Cos($lat_rad)
    $lat_tan = $lat_sin / $lat_cos
    $lat_tan2 = $lat_tan * $lat_tan
    $lat_tan4 = $lat_tan2 * $lat_tan2

    $zone_number = 48
    $central_lon = ($zone_number - 1) * 6 - 180 + 3
    $lon_rad = Radians $lon
    $central_lon_rad = Radians $central_lon

    $n = $R / $wf.Sqrt(1 - $E * $lat_sin * $lat_sin)
    $c = $E_P2 * $lat_cos * $lat_cos

    $a = $lat_cos * ($lon_rad - $central_lon_rad)
    $a2 = $a * $a
    $a3 = $a2 * $a
    $a4 = $a3 * $a
    $a5 = $a4 * $a
    $a6 = $a5 * $a

    $m = $R * ($M1 * $lat_rad - $M2 * $wf.Sin(2 * $lat_rad) + $M3 * $wf.Sin(4 * $lat_rad) - $M4 * $wf.Sin(6 * $lat_rad))

    $easting = $K0 * $n * ($a + $a3 / 6 * (1 - $lat_tan2 + $c) + $a5 / 120 * (5 - 18 * $lat_tan2 + $lat_tan4 + 72 * $c - 58 * $E_P2)) + 500000

    $northing = $K0 * ($m + $n * $lat_tan * ($a2 / 2 + $a4 / 24 * (5 - $lat_tan2 + 9 * $c + 4 * $c * $c) + $a6 / 720 * (61 - 58 * $lat_tan2 + $lat_tan4 + 600 * $c - 330 * $E_P2)))

    if ($lat -lt 0) {
        $northing = $northing + 10000000
    }

    return @($easting, $northing)
}

$firstRow = 2
$lastRow = 180

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $lat = $ws.Cells.Item($r, 9).Value()    # column I = Latitude
    $lon = $ws.Cells.Item($r, 10).Value()   # column J = Longitude

    $coords = ToUtmZone48T $lat $lon

    $ws.Cells.Item($r, 2).Value = $coords[0]   # B = Y_UTM (Easting)
    $ws.Cells.Item($r, 3).Value = $coords[1]   # C = X_UTM (Northing)
    $ws.Cells.Item($r, 4).Value = 48           # D = ZoneNumber
    $ws.Cells.Item($r, 5).Value = "T"          # E = ZoneLetter
}
